# Apply the "checking in on thurs 14-nov" edit to the TC04 worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC04")
$ws.Activate()

# Values entered in the same order the author typed them, so that the
# shared-string table is rebuilt in the same sequence as the target file.
$ws.Range("E2").Value = "Zoho CRM - Create Task"
$ws.Range("E1").Value = "tasktitle"
$ws.Range("G2").Value = "Highest"
$ws.Range("F1").Value = "Subject"
$ws.Range("G1").Value = "Priority"
$ws.Range("F2").Value = "NEW TASK CREATED"
$ws.Range("H1").Value = "lead/Contact"
$ws.Range("H2").Value = "Leads"

# Column E width adjustment (best-fit width of 22 characters for the new content)
$ws.Columns.Item(5).ColumnWidth = 21.1666666666667

# Update selection to H4 (matches the sheet's saved cursor position)
$ws.Range("H4").Select()
